$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D2:D51 as Text before writing, so numeric-looking price strings
# (e.g. "7.60", "62.728.50") are preserved verbatim as text, matching the
# source workbook which stores them as inline strings. ClearFormats()
# afterwards drops the explicit style again so cells end up unstyled,
# same as before the edit.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.728.50"
$ws.Range("E2").Value = "  +1.92%  "

$ws.Range("D3").Value = "3.466.56"
$ws.Range("E3").Value = "  +2.23%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "578.42"
$ws.Range("E5").Value = "  +0.39%  "

$ws.Range("D6").Value = "147.32"
$ws.Range("E6").Value = "  +3.89%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "0.481"
$ws.Range("E8").Value = "  +1.54%  "

$ws.Range("D9").Value = "7.60"
$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").Value = "0.124"
$ws.Range("E10").Value = "  +1.67%  "

$ws.Range("E11").Value = "  +3.88%  "

$ws.Range("D12").Value = "4.061.93"
$ws.Range("E12").Value = "  +2.37%  "

$ws.Range("D13").Value = "29.75"
$ws.Range("E13").Value = "  +5.19%  "

$ws.Range("D14").Value = "0.128"
$ws.Range("E14").Value = "  +2.42%  "

$ws.Range("D15").Value = "3.466.20"
$ws.Range("E15").Value = "  +2.52%  "

$ws.Range("E16").Value = "  +0.69%  "

$ws.Range("D17").Value = "62.840.53"
$ws.Range("E17").Value = "  +2.06%  "

$ws.Range("D18").Value = "6.34"
$ws.Range("E18").Value = "  +3.24%  "

$ws.Range("D19").Value = "14.32"
$ws.Range("E19").Value = "  +5.25%  "

$ws.Range("E20").Value = "  +2.55%  "

$ws.Range("D21").Value = "387.92"
$ws.Range("E21").Value = "  -0.73%  "

$ws.Range("D22").Value = "0.560"
$ws.Range("E22").Value = "  +2.10%  "

$ws.Range("D23").Value = "74.57"
$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "3.611.23"
$ws.Range("E25").Value = "  +2.50%  "

$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").Value = "  +1.70%  "

$ws.Range("D27").Value = "0.178"
$ws.Range("E27").Value = "  -9.78%  "

$ws.Range("D28").Value = "7.57"
$ws.Range("E28").Value = "  +2.95%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").Value = "8.16"
$ws.Range("E30").Value = "  +2.13%  "

$ws.Range("D31").Value = "2.13"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("E33").Value = "  -1.00%  "

$ws.Range("D34").Value = "23.70"
$ws.Range("E34").Value = "  +1.79%  "

$ws.Range("D35").Value = "5.28"
$ws.Range("E35").Value = "  +5.04%  "

$ws.Range("D36").Value = "7.08"
$ws.Range("E36").Value = "  +2.60%  "

$ws.Range("D37").Value = "1.59"
$ws.Range("E37").Value = "  +7.91%  "

$ws.Range("D38").Value = "31.70"
$ws.Range("E38").Value = "  +22.84%  "

$ws.Range("D39").Value = "170.03"
$ws.Range("E39").Value = "  +0.48%  "

$ws.Range("D40").Value = "3.505.72"
$ws.Range("E40").Value = "  +2.44%  "

$ws.Range("D41").Value = "0.0765"
$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").Value = "0.799"
$ws.Range("E42").Value = "  +2.62%  "

$ws.Range("D43").Value = "4.47"
$ws.Range("E43").Value = "  +0.97%  "

$ws.Range("D44").Value = "42.22"
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("D45").Value = "1.71"
$ws.Range("E45").Value = "  +4.21%  "

$ws.Range("E46").Value = "  +3.68%  "

$ws.Range("D47").Value = "2.609.92"
$ws.Range("E47").Value = "  +5.88%  "

$ws.Range("D48").Value = "23.13"
$ws.Range("E48").Value = "  +1.66%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  +10.77%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "6.74"
$ws.Range("E50").Value = "  +1.29%  "

$ws.Range("E51").Value = "  +0.18%  "

$priceRange.ClearFormats()
